$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.723.13'
$ws.Range('D3').Value = '2.620.00'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.94'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.46'
$ws.Range('E6').Value = '  +2.59%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.589'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.110'
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.70'
$ws.Range('E10').Value = '  +2.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.382'
$ws.Range('E11').Value = '  +3.20%  '
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.76'
$ws.Range('E13').Value = '  +0.44%  '
$ws.Range('D14').Value = '3.090.97'
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('D15').Value = '63.493.33'
$ws.Range('E15').Value = '  +0.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000155'
$ws.Range('E16').Value = '  +5.10%  '
$ws.Range('D17').Value = '2.621.14'
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.33'
$ws.Range('E18').Value = '  +7.04%  '
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '347.17'
$ws.Range('E20').Value = '  +0.72%  '
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.73'
$ws.Range('E23').Value = '  +2.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.36'
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('E25').Value = '  +11.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.28'
$ws.Range('E26').Value = '  +1.91%  '
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '563.86'
$ws.Range('E28').Value = '  -3.29%  '
$ws.Range('E29').Value = '  +3.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.162'
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.05'
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '168.54'
$ws.Range('E36').Value = '  +0.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.409'
$ws.Range('E37').Value = '  +0.42%  '
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.38'
$ws.Range('E40').Value = '  +1.32%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '167.05'
$ws.Range('E42').Value = '  -0.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.89'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('E44').Value = '  +3.42%  '
$ws.Range('E45').Value = '  +4.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.55'
$ws.Range('E46').Value = '  -3.00%  '
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('E49').Value = '  +4.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0964'
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.37'
$ws.Range('E51').Value = '  +3.07%  '
